$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving format for numeric-looking price strings so Excel
# does not auto-convert them to floating point numbers.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D35", "D38", "D39", "D40", "D41", "D42", "D46", "D47", "D49", "D50", "D51")
foreach ($cell in $textCells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range("D2").Value = "35.312.12"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.911.94"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "0.719"
$ws.Range("E5").Value = "  +9.19%  "
$ws.Range("D6").Value = "253.79"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "40.73"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "0.360"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("D10").Value = "52.27"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +4.23%  "
$ws.Range("D12").Value = "0.0992"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "2.188.88"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "12.61"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("D15").Value = "0.718"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "4.91"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "1.903.41"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "35.324.06"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "74.14"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "243.62"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "12.98"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").Value = "5.08"
$ws.Range("E23").Value = "  +3.96%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +5.52%  "
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").Value = "166.82"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "18.72"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +4.03%  "
$ws.Range("D31").Value = "4.130.19"
$ws.Range("E31").Value = "  +19.50%  "
$ws.Range("D32").Value = "4.35"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("D33").Value = "2.00"
$ws.Range("E33").Value = "  +14.78%  "
$ws.Range("E34").Value = "  +22.17%  "
$ws.Range("D35").Value = "0.0581"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "0.915"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "2.03"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "17.34"
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("D41").Value = "0.0217"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").Value = "97.01"
$ws.Range("E42").Value = "  +6.86%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "1.338.03"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "2.77"
$ws.Range("E49").Value = "  -1.39%  "

# Row 50/51: Gas and MultiversX swapped rank position in this update
$ws.Range("B50").Value = "Gas"
$ws.Range("C50").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D50").Value = "12.36"
$ws.Range("E50").Value = "  +22.50%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "45.14"
$ws.Range("E51").Value = "  -7.33%  "

# Restore the default (unstyled) cell style now that the text values are set,
# matching the original workbook formatting.
foreach ($cell in $textCells) { $ws.Range($cell).Style = "Normal" }

Write-Output "done"
